# updating meta email id
#
# The Kawasaki Daily Reporting Meta datorama report id changed
# (1159468 -> 1161493). Both the NAV sheet (which already had the old id)
# and the TeryxH2 sheet (whose meta_gmail_subject cell was still blank)
# need to show the new subject line. The previously-blank TeryxH2 cell
# also loses its "needs attention" yellow highlight now that it's filled
# in. Finally, the workbook's active tab moves from TeryxH2 back to NAV.

$wb = $excel.ActiveWorkbook

$newSubject = 'Datorama | Report "Kawasaki Daily Reporting Meta" (1161493)'

# --- NAV sheet: meta_gmail_subject (B9) gets the refreshed report id ---
$nav = $wb.Worksheets.Item("NAV")
$nav.Range("B9").Value = $newSubject

# --- TeryxH2 sheet: meta_gmail_subject (B9) was blank, now filled in,   ---
# --- and the yellow "missing" highlight fill is cleared.                ---
$teryx = $wb.Worksheets.Item("TeryxH2")
$teryx.Range("B9").Value = $newSubject
$teryx.Range("B9").Interior.Pattern = -4142
$teryx.Range("B9").Interior.ColorIndex = -4142

# --- View state: NAV becomes the active/selected tab with B9 selected, ---
# --- TeryxH2 is no longer the active tab and its selection moves to    ---
# --- B10 (first empty row below the table). Set TeryxH2's selection    ---
# --- first, then activate NAV last so NAV ends up as the visible tab.  ---
$teryx.Range("B10").Select()

$nav.Activate()
$nav.Range("B9").Select()
